# "Use namePhrase instead of appendedPhrase"
#
# The workbook had two near-duplicate shared strings:
#   - "appendedNamePhrase" used as the header of column D on the Taxon sheet
#   - "appendedPhrase"     used as the header of column C on the Synonym sheet
# Both are replaced by a single shared header "namePhrase", and the Taxon
# sheet's new namePhrase header cell picks up the same (Helvetica) font
# already used for the namePhrase header on the Synonym sheet.
#
# The previously active sheet/tab (TaxonRelation) is swapped for Taxon,
# with Taxon's header cell (D1) and Synonym's header cell (C1) selected.

$wb = $excel.ActiveWorkbook

$taxon    = $wb.Worksheets.Item("Taxon")
$synonym  = $wb.Worksheets.Item("Synonym")

# Rename the header cells to the unified "namePhrase" label.
$taxon.Range("D1").Value = "namePhrase"
$synonym.Range("C1").Value = "namePhrase"

# Match the formatting already used for the "namePhrase"/"appendedPhrase"
# header on the Synonym sheet (Helvetica font).
$taxon.Range("D1").Font.Name = "Helvetica"

# Update each sheet's remembered selection.
$synonym.Activate()
$synonym.Range("C1").Select()

# Make Taxon the active sheet/tab, selecting its renamed header cell.
$taxon.Activate()
$taxon.Range("D1").Select()
